$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 8 of trade data
$ws.Range("A8").Value = 9794.24
$ws.Range("B8").Value = 9835.5499999999993
$ws.Range("C8").Value = 80.11
$ws.Range("D8").Value = 79.77
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = -0.42

# G8 holds a date/time value, styled like the other date cells in column G (style index 1 -> numFmtId 22)
$ws.Range("G8").Value = 42609.488321759258
$ws.Range("G8").NumberFormat = "m/d/yy h:mm"

$ws.Range("H8").Value = $false
